$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("C2").Value = "sadasd"
$ws.Range("D2").Value = "sdasd"
$ws.Range("F2").Value = "12:10"

# Row 3 updates
$ws.Range("C3").Value = "sadasd"
$ws.Range("D3").Value = "sdasd"

# Remove publish_date / publish_time for row 3 entirely (clear the cells)
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
